# Applies the "Deduce" skill / character-board update described in the diff:
#  - Actions sheet: insert a new row 22 ("Deduce" / "💡💡🔊" / "Ideas"),
#    pushing the existing rows 22-46 down to 23-47.
#  - Characters sheet: update the L8/M8 flavor-text formulas, and shuffle the
#    G9/G10/H10 action labels (Report -> Deduce at G9, Inspire -> Report at
#    G10, and Inspire now also appears at H10).
#  - Selections on both sheets move to reflect the edited cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Actions sheet - insert the new "Deduce" row above the current row 22
# ---------------------------------------------------------------------
$actions = $wb.Worksheets.Item("Actions")

$actions.Rows("22:22").Insert()

$actions.Range("A22").Value = "Deduce"
$actions.Range("B22").Value = "💡💡🔊"
$actions.Range("C22").Value = "Ideas"
$actions.Range("D22").Formula = "=COUNTIF(Skills!E:E,A22) + COUNTIF(Skills!I:I,A22) + COUNTIF(Skills!K:K,A22) + COUNTIF(Skills!M:M,A22) + COUNTIF(Skills!O:O,A22)"

[void]$actions.Range("A22:B22").Select()

# The "Actions" defined name covers the data rows (A2:B<lastRow>) and must
# grow by one row along with the freshly-inserted row.
$wb.Names.Item("Actions").RefersTo = "=Actions!`$A`$2:`$B`$48"

# ---------------------------------------------------------------------
# Characters sheet - flavor text + action-label updates
# ---------------------------------------------------------------------
$characters = $wb.Worksheets.Item("Characters")

$characters.Range("L8").Formula = '="+3 Memory%n+1 Initial 💡%n %nWalk🔊➜ ⇒%nDash🔊➜➜%n %nReport💡🔍🔊⇒%nDeduce💡💡🔊"'
$characters.Range("M8").Formula = '="+1 Memory%n+1 Initial 💡%n %nWalk🔊➜ ⇒%nRun🔊🔊➜➜%n %n+INSPIRE"'

$characters.Range("G9").Value = "Deduce%n💡💡🔊"
$characters.Range("G10").Value = "Report%n💡🔍🔊"
$characters.Range("H10").Value = "Inspire"

[void]$characters.Range("K9").Select()
